{"js": "// Update the date heading (first paragraph) and regenerate the\n// division-fact table contents for the new day's worksheet.\nconst body = context.document.body;\n\n// 1) Update the title paragraph: \"2024-05-21 Tuesday\" -> \"2024-05-22 Wednesday\"\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst titlePara = paragraphs.items[0];\ntitlePara.load(\"text\");\nawait context.sync();\n\nif (titlePara.text.indexOf(\"2024-05-21 Tuesday\") !== -1) {\n  titlePara.insertText(\"2024-05-22 Wednesday\", \"Replace\");\n}\n\n// 2) Update the practice table. The table alternates a row of five\n// division-fact cells with three blank spacer rows, five times over.\n// Replace the five-cell \"fact\" rows' text in place (left-to-right,\n// top-to-bottom) while leaving the blank spacer rows untouched.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst newFacts = [\n  \"72\u00f79=\", \"56\u00f72=\", \"13\u00f76=\", \"77\u00f78=\", \"66\u00f74=\",\n  \"41\u00f75=\", \"73\u00f75=\", \"38\u00f72=\", \"94\u00f73=\", \"71\u00f73=\",\n  \"53\u00f76=\", \"57\u00f76=\", \"44\u00f77=\", \"24\u00f79=\", \"85\u00f79=\",\n  \"76\u00f74=\", \"83\u00f74=\", \"73\u00f74=\", \"80\u00f78=\", \"41\u00f74=\",\n  \"63\u00f79=\", \"45\u00f79=\", \"52\u00f73=\", \"22\u00f77=\", \"65\u00f72=\",\n];\n\nconst oldValues = table.values;\nlet factIndex = 0;\nconst newValues = oldValues.map((row) => {\n  const hasText = row.some((cell) => cell !== \"\");\n  if (!hasText) {\n    // Blank spacer row: leave as-is.\n    return row;\n  }\n  return row.map(() => newFacts[factIndex++]);\n});\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# Update the date heading and regenerate the division-fact table\n# contents for the new day's worksheet.\n$d = $word.ActiveDocument\n\n# 1) Update the title paragraph: \"2024-05-21 Tuesday\" -> \"2024-05-22 Wednesday\"\n$titlePara = $d.Paragraphs.Item(1)\nif ($titlePara.Range.Text -like \"*2024-05-21 Tuesday*\") {\n    $titlePara.Range.Text = \"2024-05-22 Wednesday\"\n}\n\n# 2) Update the practice table. The table alternates a row of five\n# division-fact cells with three blank spacer rows, five times over.\n# Replace the five-cell \"fact\" rows' text in place (left-to-right,\n# top-to-bottom, by cell index) while leaving the blank spacer rows\n# untouched.\n$newFacts = @(\n    \"72\u00f79=\", \"56\u00f72=\", \"13\u00f76=\", \"77\u00f78=\", \"66\u00f74=\",\n    \"41\u00f75=\", \"73\u00f75=\", \"38\u00f72=\", \"94\u00f73=\", \"71\u00f73=\",\n    \"53\u00f76=\", \"57\u00f76=\", \"44\u00f77=\", \"24\u00f79=\", \"85\u00f79=\",\n    \"76\u00f74=\", \"83\u00f74=\", \"73\u00f74=\", \"80\u00f78=\", \"41\u00f74=\",\n    \"63\u00f79=\", \"45\u00f79=\", \"52\u00f73=\", \"22\u00f77=\", \"65\u00f72=\"\n)\n\n$table = $d.Tables.Item(1)\n$factRows = @(1, 5, 9, 13, 17)\n\n$factIndex = 0\nforeach ($r in $factRows) {\n    for ($c = 1; $c -le 5; $c++) {\n        $table.Cell($r, $c).Range.Text = $newFacts[$factIndex]\n        $factIndex++\n    }\n}\n"}
